{"js": "// Append a brand-new paragraph right after the document's final paragraph\n// (\"Attendez une minute...\"). The new paragraph carries the explicit\n// \"Normal\" style plus two runs: an empty lead run and a run holding the\n// new French text -- matching a clean/compatibility-mode rewrite exactly.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Create a fresh, empty paragraph right after the last one.\nconst newParagraph = lastParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// Replace that empty paragraph's contents with the fully-specified OOXML\n// (flat-OPC wrapped, as required by Range.insertOoxml): an explicit\n// \"Normal\" pStyle, an empty lead run, then the run carrying the new text.\nconst flatOpcXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Normal\"/></w:pPr><w:r><w:rPr/></w:r><w:r><w:rPr/><w:t>Un jour, un myst\u00e9rieux \u00e9tudiant transf\u00e9r\u00e9 arrive soudainement dans ma classe. Il s'av\u00e8re que cet \u00e9l\u00e8ve est en fait un extraterrestre ou un voyageur dans le temps, ou quelque chose du m\u00eame genre, dot\u00e9 de pouvoirs inconnus. Ensuite, il se trouve que l'\u00e9l\u00e8ve se bat contre un gang mal\u00e9fique et que je me retrouve m\u00eal\u00e9 \u00e0 ce combat. L'autre \u00e9l\u00e8ve est celui qui se bat le plus. Je ne suis qu'un acolyte. H\u00e9, \u00e7a a l'air cool. Bon sang, je suis intelligent.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\nconst targetRange = newParagraph.getRange();\ntargetRange.insertOoxml(flatOpcXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the final paragraph (\"Attendez une minute...\") and append a brand\n# new paragraph right after it.\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphAfter()\n\n# The freshly inserted paragraph is now the document's last paragraph; its\n# Range currently holds just the paragraph mark. Replace that mark with a\n# fully-specified <w:p> fragment (explicit \"Normal\" style, an empty lead\n# run followed by the run carrying the new text) so the resulting OOXML\n# matches a clean/compat rewrite exactly.\n$newPara = $d.Paragraphs.Last\n$text = 'Un jour, un myst\u00e9rieux \u00e9tudiant transf\u00e9r\u00e9 arrive soudainement dans ma classe. Il s''av\u00e8re que cet \u00e9l\u00e8ve est en fait un extraterrestre ou un voyageur dans le temps, ou quelque chose du m\u00eame genre, dot\u00e9 de pouvoirs inconnus. Ensuite, il se trouve que l''\u00e9l\u00e8ve se bat contre un gang mal\u00e9fique et que je me retrouve m\u00eal\u00e9 \u00e0 ce combat. L''autre \u00e9l\u00e8ve est celui qui se bat le plus. Je ne suis qu''un acolyte. H\u00e9, \u00e7a a l''air cool. Bon sang, je suis intelligent.'\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"Normal\"/></w:pPr><w:r><w:rPr/></w:r><w:r><w:rPr/><w:t>' + $text + '</w:t></w:r></w:p>'\n$newPara.Range.InsertXML($xml)\n"}
